$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$abstractText = @"
Objectives
Rapid and accurate diagnosis of coronavirus disease 2019 (COVID-19) is critical during the epidemic.
 We aim to identify differences in CT imaging and clinical manifestations between pneumonia patients with and without COVID-19, and to develop and validate a diagnostic model for COVID-19 based on radiological semantic and clinical features alone.
Methods
A consecutive cohort of 70 COVID-19 and 66 non-COVID-19 pneumonia patients were retrospectively recruited from five institutions.
 Patients were divided into primary (n = 98) and validation (n = 38) cohorts.
 The chi-square test, Student’s t test, and Kruskal-Wallis H test were performed, comparing 1745 lesions and 67 features in the two groups.
 Three models were constructed using radiological semantic and clinical features through multivariate logistic regression.
 Diagnostic efficacies of developed models were quantified by receiver operating characteristic curve.
 Clinical usage was evaluated by decision curve analysis and nomogram.
Results
Eighteen radiological semantic features and seventeen clinical features were identified to be significantly different.
 Besides ground-glass opacities (p = 0.032) and consolidation (p = 0.001) in the lung periphery, the lesion size (1–3 cm) is also significant for the diagnosis of COVID-19 (p = 0.027).
 Lung score presents no significant difference (p = 0.417).
 Three diagnostic models achieved an area under the curve value as high as 0.986 (95% CI 0.966~1.000).
 The clinical and radiological semantic models provided a better diagnostic performance and more considerable net benefits.
Conclusions
Based on CT imaging and clinical manifestations alone, the pneumonia patients with and without COVID-19 can be distinguished.
 A model composed of radiological semantic and clinical features has an excellent performance for the diagnosis of COVID-19.
Key Points
            • Based on CT imaging and clinical manifestations alone, the pneumonia patients with and without COVID-19 can be distinguished.
Electronic supplementary material
The online version of this article (10.1007/s00330-020-06829-2) contains supplementary material, which is available to authorized users.

"@

$authorsText = @"
[Xiaofeng%Chen%NULL%0,      Yanyan%Tang%NULL%1,      Yongkang%Mo%NULL%1,      Shengkai%Li%NULL%1,      Daiying%Lin%NULL%1,      Zhijian%Yang%NULL%1,      Zhiqi%Yang%NULL%1,      Hongfu%Sun%NULL%1,      Jinming%Qiu%NULL%1,      Yuting%Liao%NULL%1,      Jianning%Xiao%NULL%1,      Xiangguang%Chen%NULL%1,      Xianheng%Wu%NULL%1,      Renhua%Wu%NULL%1,      Zhuozhi%Dai%zhuozhi@ualberta.ca%1]
"@

$ws.Range("D6").Value = $abstractText
$ws.Range("E6").Value = $authorsText
